# Applies the diff:
#  1. Removes the _GoBack bookmark that sat right after the "precision metric"
#     sentence.
#  2. Fills the empty numbered-list paragraph ("Final Agent Performance" bullet
#     child) with the "I can see from the graph..." sentence, re-inserting a
#     _GoBack bookmark mid-sentence (after the hyphen in "net-rewad").
#  3. Turns the next (plain, spacing-only) empty paragraph into a numbered
#     list item too (same list/level as its sibling) and fills it with the
#     "For some reason..." / zero-crossing sentence.

$d = $word.ActiveDocument

# --- Helper fragments -------------------------------------------------
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPrTNR = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

# --- Step 1: drop the old _GoBack bookmark -----------------------------
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
    # no-op if it is already gone
}

# --- Step 2: find the two target paragraphs ----------------------------
# Paragraph A: empty numbered ("ListParagraph", ilvl=1, numId=6) bullet
#              right after "Final Agent Performance" / "Referring ..." item.
# Paragraph B: the following empty paragraph (spacing-only pPr, no list).
$paraA = $null
$paraB = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -eq 0 -and $p.Range.ListFormat.ListType -ne 0) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim().Length -eq 0 -and $next.Range.ListFormat.ListType -eq 0) {
            $paraA = $p
            $paraB = $next
            break
        }
    }
}

if ($paraA -eq $null -or $paraB -eq $null) {
    throw "Could not locate the target empty paragraphs"
}

# --- Step 3: fill paragraph A with its runs + the new _GoBack bookmark --
$bodyA = '<w:p>' +
    '<w:r>' + $rPrTNR + '<w:t xml:space="preserve">I can see from the </w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t>graph I drew at the end of 100 trials that the agent tends to keep the net</w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t>-</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r>' + $rPrTNR + '<w:t>rewad as positive consistently indicating me that it has learned and very close to the stated optimal policy</w:t></w:r>' +
    '</w:p>'
$xmlA = $pkgOpen + $bodyA + $pkgClose

$rangeA = $paraA.Range
$insertAtA = $d.Range($rangeA.End - 1, $rangeA.End - 1)
$null = $insertAtA.InsertXML($xmlA)

# --- Step 4: rebuild paragraph B (pPr + runs) in one shot ---------------
$bodyB = '<w:p>' +
    '<w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr>' +
      '<w:spacing w:after="225" w:line="357" w:lineRule="atLeast"/>' +
      '<w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Times New Roman"/><w:color w:val="58646D"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' + $rPrTNR + '<w:t xml:space="preserve">For some reason that I am not certain, the net reward hasn&#8217;t gone below zero at any time which makes me wonder whether it was a loose implementation of the reward rule or whether I am doing something wrong. </w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t>If there</w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t xml:space="preserve"> was a </w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t>zero</w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t>-crossing as it described on the referenced site, then it will tell us how long the algorithm or policy takes to recoup the cost of learning.</w:t></w:r>' +
    '<w:r>' + $rPrTNR + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$xmlB = $pkgOpen + $bodyB + $pkgClose

# paraB reference is still valid (we only edited paraA, which sits before it)
$rangeB = $paraB.Range
$fullB = $d.Range($rangeB.Start, $rangeB.End)
$null = $fullB.InsertXML($xmlB)

Write-Host "Done"
